$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update relocated/renamed suggestion values in column B
$ws.Range("B9").Value = "ギャラクシー原宿"
$ws.Range("B11").Value = "ギャラクシーs10 カバー"

# Remove the now-unused rows (12-31) that held the duplicate Galaxy rows
# and the RPA keyword block
$ws.Range("A12:B31").EntireRow.Delete()
